# Remove the "exception screenshots" settings row (and the blank spacer row
# beneath it) from the Constants sheet. Everything below shifts up by two
# rows as a result, matching the rest of the sheet's layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

$ws.Rows("5:6").Delete()

$ws.Range("B11").Select() | Out-Null
